# Auto-generated script to apply numeric corrections to Atomos_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3139.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3139.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3139.8
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -3277.8
$ws.Range("H95").Value = 22887.3
$ws.Range("J95").Value = 22887.3
$ws.Range("L95").Value = 22887.3
$ws.Range("N95").Value = -28379.3
$ws.Range("H116").Value = 3609.4482
$ws.Range("I116").Value = 2886.7058
$ws.Range("K116").Value = 2886.7058
$ws.Range("M116").Value = 555.2941999999998
$ws.Range("H132").Value = 6670603
$ws.Range("I132").Value = 8003107.5
$ws.Range("K132").Value = 24009322.5
$ws.Range("M132").Value = -24006792.5
$ws.Range("H136").Value = 29000
$ws.Range("J136").Value = 29000
$ws.Range("L136").Value = 29000
$ws.Range("N136").Value = -39200
$ws.Range("H137").Value = 3988.525
$ws.Range("I137").Value = 4655.2
$ws.Range("J137").Value = 2877.4
$ws.Range("K137").Value = 13965.6
$ws.Range("L137").Value = 8632.200000000001
$ws.Range("M137").Value = -11415.6
$ws.Range("N137").Value = -13732.2
$ws.Range("H141").Value = 393448.3
$ws.Range("I141").Value = 998.2632
$ws.Range("J141").Value = 807701.1
$ws.Range("K141").Value = 2994.7896
$ws.Range("L141").Value = 2423103.3
$ws.Range("M141").Value = 2185.2104
$ws.Range("N141").Value = -2433463.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2943.4375
$ws.Range("I74").Value = 2009.2
$ws.Range("K74").Value = 2009.2
$ws.Range("M74").Value = -1135.2
$ws.Range("H77").Value = 2943.4375
$ws.Range("I77").Value = 2009.2
$ws.Range("K77").Value = 10046
$ws.Range("M77").Value = -5678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 30000
$ws.Range("J132").Value = 30000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -40120
$ws.Range("H134").Value = 6226.6
$ws.Range("I134").Value = 5898.478
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 17695.434
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -15160.434
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 34091.363
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 37000.5
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 37000.5
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -37224.5
$ws.Range("H58").Value = 20838860
$ws.Range("I58").Value = 4293.091
$ws.Range("J58").Value = 38468108
$ws.Range("K58").Value = 4293.091
$ws.Range("L58").Value = 38468108
$ws.Range("M58").Value = -4090.091
$ws.Range("N58").Value = -38468514
$ws.Range("H105").Value = 3654.4443
$ws.Range("I105").Value = 7003.3335
$ws.Range("J105").Value = 1980
$ws.Range("K105").Value = 7003.3335
$ws.Range("L105").Value = 1980
$ws.Range("M105").Value = -5256.3335
$ws.Range("N105").Value = -5474
$ws.Range("H136").Value = 20838860
$ws.Range("I136").Value = 4293.091
$ws.Range("J136").Value = 38468108
$ws.Range("K136").Value = 12879.273
$ws.Range("L136").Value = 115404324
$ws.Range("M136").Value = -10329.273
$ws.Range("N136").Value = -115409424
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 364210.62
$ws.Range("I4").Value = 800062.6
$ws.Range("J4").Value = 1000.6667
$ws.Range("K4").Value = 2400187.8
$ws.Range("L4").Value = 3002.0001
$ws.Range("M4").Value = -2400075.8
$ws.Range("N4").Value = -3226.0001
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9920
$ws.Range("H60").Value = 7326.316
$ws.Range("I60").Value = 200
$ws.Range("J60").Value = 8662.5
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 25987.5
$ws.Range("M60").Value = -349
$ws.Range("N60").Value = -26489.5
$ws.Range("H87").Value = 14457.143
$ws.Range("I87").Value = 10950
$ws.Range("J87").Value = 15860
$ws.Range("K87").Value = 32850
$ws.Range("L87").Value = 47580
$ws.Range("M87").Value = -31602
$ws.Range("N87").Value = -50076
$ws.Range("H90").Value = 14457.143
$ws.Range("I90").Value = 10950
$ws.Range("J90").Value = 15860
$ws.Range("K90").Value = 98550
$ws.Range("L90").Value = 142740
$ws.Range("M90").Value = -92310
$ws.Range("N90").Value = -155220
$ws.Range("H94").Value = 3762.3076
$ws.Range("J94").Value = 3992.1738
$ws.Range("L94").Value = 11976.5214
$ws.Range("N94").Value = -13328.5214
$ws.Range("H109").Value = 898.9375
$ws.Range("I109").Value = 298.3
$ws.Range("J109").Value = 1900
$ws.Range("K109").Value = 894.9000000000001
$ws.Range("L109").Value = 5700
$ws.Range("M109").Value = 145.0999999999999
$ws.Range("N109").Value = -7780
$ws.Range("H134").Value = 2617.7856
$ws.Range("I134").Value = 1853.2222
$ws.Range("J134").Value = 3994
$ws.Range("K134").Value = 5559.6666
$ws.Range("L134").Value = 11982
$ws.Range("M134").Value = -489.6665999999996
$ws.Range("N134").Value = -22122
$ws.Range("H136").Value = 2449.0557
$ws.Range("I136").Value = 1686.6666
$ws.Range("J136").Value = 3211.4443
$ws.Range("K136").Value = 5059.9998
$ws.Range("L136").Value = 9634.332900000001
$ws.Range("M136").Value = 40.0002000000004
$ws.Range("N136").Value = -19834.3329
$ws.Range("H137").Value = 2818.95
$ws.Range("I137").Value = 2098.5
$ws.Range("J137").Value = 4500
$ws.Range("K137").Value = 6295.5
$ws.Range("L137").Value = 13500
$ws.Range("M137").Value = -1195.5
$ws.Range("N137").Value = -23700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12333.333
$ws.Range("I5").Value = 8500
$ws.Range("J5").Value = 20000
$ws.Range("K5").Value = 8500
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = -8388
$ws.Range("N5").Value = -20224
$ws.Range("H18").Value = 190167
$ws.Range("I18").Value = 500650
$ws.Range("J18").Value = 86672.664
$ws.Range("K18").Value = 500650
$ws.Range("L18").Value = 86672.664
$ws.Range("M18").Value = -500357
$ws.Range("N18").Value = -87258.664
$ws.Range("H43").Value = 3057.0908
$ws.Range("I43").Value = 1398.3334
$ws.Range("J43").Value = 5047.6
$ws.Range("K43").Value = 1398.3334
$ws.Range("L43").Value = 5047.6
$ws.Range("M43").Value = -1247.3334
$ws.Range("N43").Value = -5349.6
$ws.Range("H46").Value = 12023
$ws.Range("J46").Value = 15046
$ws.Range("L46").Value = 15046
$ws.Range("N46").Value = -15358
$ws.Range("H57").Value = 14590.25
$ws.Range("J57").Value = 16180.5
$ws.Range("L57").Value = 16180.5
$ws.Range("N57").Value = -17820.5
$ws.Range("H80").Value = 3083.8845
$ws.Range("I80").Value = 2742.647
$ws.Range("J80").Value = 3728.4443
$ws.Range("K80").Value = 2742.647
$ws.Range("L80").Value = 3728.4443
$ws.Range("M80").Value = -1744.647
$ws.Range("N80").Value = -5724.4443
$ws.Range("H83").Value = 3083.8845
$ws.Range("I83").Value = 2742.647
$ws.Range("J83").Value = 3728.4443
$ws.Range("K83").Value = 13713.235
$ws.Range("L83").Value = 18642.2215
$ws.Range("M83").Value = -8721.235000000001
$ws.Range("N83").Value = -28626.2215
$ws.Range("H132").Value = 3279.4285
$ws.Range("I132").Value = 3119.5483
$ws.Range("J132").Value = 3554.7778
$ws.Range("K132").Value = 9358.644899999999
$ws.Range("L132").Value = 10664.3334
$ws.Range("M132").Value = -6828.644899999999
$ws.Range("N132").Value = -15724.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2363.625
$ws.Range("I7").Value = 1832.1538
$ws.Range("K7").Value = 1832.1538
$ws.Range("M7").Value = -1720.1538
$ws.Range("H22").Value = 565.2727
$ws.Range("I22").Value = 441.5
$ws.Range("J22").Value = 636
$ws.Range("K22").Value = 441.5
$ws.Range("L22").Value = 636
$ws.Range("M22").Value = -146.5
$ws.Range("N22").Value = -1226
$ws.Range("H24").Value = 100007
$ws.Range("J24").Value = 100007
$ws.Range("L24").Value = 100007
$ws.Range("N24").Value = -100693
$ws.Range("H27").Value = 565.2727
$ws.Range("I27").Value = 441.5
$ws.Range("J27").Value = 636
$ws.Range("K27").Value = 441.5
$ws.Range("L27").Value = 636
$ws.Range("M27").Value = -334.5
$ws.Range("N27").Value = -850
$ws.Range("H46").Value = 2410
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 2600
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 2600
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -2976
$ws.Range("H61").Value = 125005590
$ws.Range("I61").Value = 200002560
$ws.Range("J61").Value = 10665
$ws.Range("K61").Value = 200002560
$ws.Range("L61").Value = 10665
$ws.Range("M61").Value = -200002358
$ws.Range("N61").Value = -11069
$ws.Range("H103").Value = 32127.273
$ws.Range("J103").Value = 32127.273
$ws.Range("L103").Value = 32127.273
$ws.Range("N103").Value = -34471.273
$ws.Range("H113").Value = 125005590
$ws.Range("I113").Value = 200002560
$ws.Range("J113").Value = 10665
$ws.Range("K113").Value = 200002560
$ws.Range("L113").Value = 10665
$ws.Range("M113").Value = -200000390
$ws.Range("N113").Value = -15005
$ws.Range("H126").Value = 2363.625
$ws.Range("I126").Value = 1832.1538
$ws.Range("K126").Value = 5496.4614
$ws.Range("M126").Value = -3026.4614
$ws.Range("H132").Value = 2246.204
$ws.Range("I132").Value = 1564.9678
$ws.Range("J132").Value = 3419.4443
$ws.Range("K132").Value = 4694.903399999999
$ws.Range("L132").Value = 10258.3329
$ws.Range("M132").Value = -2164.903399999999
$ws.Range("N132").Value = -15318.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2561102
$ws.Range("I2").Value = 2062201.6
$ws.Range("K2").Value = 2062201.6
$ws.Range("M2").Value = -2062089.6
